$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-13 17:18:25"
$ws.Range("I2").Value = "1.2 mm"
$ws.Range("K2").Value = "3.4 MJ/m2"
$ws.Range("O2").Value = "-0.7 °C"
$ws.Range("E3").Value = "2026-02-13 17:18:27"
$ws.Range("I3").Value = "5.2 mm"
$ws.Range("E4").Value = "2026-02-13 17:18:29"
$ws.Range("I4").Value = "4.4 mm"
$ws.Range("J4").Value = "995.9 hPa"
$ws.Range("E5").Value = "2026-02-13 17:18:32"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "80%"
$ws.Range("I5").Value = "0.1 mm"
$ws.Range("K5").Value = "4.8 MJ/m2"
$ws.Range("E6").Value = "2026-02-13 17:18:34"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "73%"
$ws.Range("I6").Value = "3.0 mm"
$ws.Range("J6").Value = "995.9 hPa"
$ws.Range("E7").Value = "2026-02-13 17:18:37"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "70%"
$ws.Range("I7").Value = "15.4 mm"
$ws.Range("J7").Value = "996.1 hPa"
$ws.Range("N7").Value = "10.9 °C 16:31 TU"
$ws.Range("O7").Value = "13.1 °C"
$ws.Range("E8").Value = "2026-02-13 17:18:39"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "77%"
$ws.Range("I8").Value = "17.3 mm"
$ws.Range("J8").Value = "996.0 hPa"
$ws.Range("O8").Value = "9.6 °C"
$ws.Range("E9").Value = "2026-02-13 17:18:41"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "73%"
$ws.Range("I9").Value = "1.1 mm"
$ws.Range("E10").Value = "2026-02-13 17:18:44"
$ws.Range("I10").Value = "16.1 mm"
$ws.Range("E11").Value = "2026-02-13 17:18:46"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "91%"
$ws.Range("E12").Value = "2026-02-13 17:18:49"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "83%"
$ws.Range("I12").Value = "4.0 mm"
$ws.Range("E13").Value = "2026-02-13 17:18:51"
$ws.Range("E14").Value = "2026-02-13 17:18:54"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "82%"
$ws.Range("I14").Value = "17.7 mm"
$ws.Range("E15").Value = "2026-02-13 17:18:56"
$ws.Range("I15").Value = "1.5 mm"
$ws.Range("E16").Value = "2026-02-13 17:18:59"
$ws.Range("G16").Value = "82 cm"
$ws.Range("K16").Value = "6.2 MJ/m2"
$ws.Range("E17").Value = "2026-02-13 17:19:01"
$ws.Range("I17").Value = "4.6 mm"
$ws.Range("O17").Value = "0.6 °C"
$ws.Range("E18").Value = "2026-02-13 17:19:04"
$ws.Range("I18").Value = "6.7 mm"
$ws.Range("J18").Value = "996.1 hPa"
$ws.Range("E19").Value = "2026-02-13 17:19:06"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "88%"
$ws.Range("I19").Value = "10.8 mm"
$ws.Range("K19").Value = "2.1 MJ/m2"
$ws.Range("E20").Value = "2026-02-13 17:19:09"
$ws.Range("I20").Value = "19.1 mm"
$ws.Range("E21").Value = "2026-02-13 17:19:11"
$ws.Range("J21").Value = "998.8 hPa"
$ws.Range("K21").Value = "1.7 MJ/m2"
$ws.Range("E22").Value = "2026-02-13 17:19:13"
$ws.Range("K22").Value = "6.8 MJ/m2"
$ws.Range("E23").Value = "2026-02-13 17:19:16"
$ws.Range("I23").Value = "7.2 mm"
$ws.Range("E24").Value = "2026-02-13 17:19:18"
$ws.Range("J24").Value = "996.4 hPa"
$ws.Range("E25").Value = "2026-02-13 17:19:21"
$ws.Range("I25").Value = "8.4 mm"
$ws.Range("K25").Value = "4.0 MJ/m2"
$ws.Range("E26").Value = "2026-02-13 17:19:23"
$ws.Range("E27").Value = "2026-02-13 17:19:26"
$ws.Range("O27").Value = "-2.4 °C"
$ws.Range("E28").Value = "2026-02-13 17:19:28"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "77%"
$ws.Range("I28").Value = "6.2 mm"
$ws.Range("J28").Value = "996.4 hPa"
$ws.Range("O28").Value = "6.5 °C"
$ws.Range("E29").Value = "2026-02-13 17:19:31"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "85%"
$ws.Range("I29").Value = "11.6 mm"
$ws.Range("E30").Value = "2026-02-13 17:19:33"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "76%"
$ws.Range("I30").Value = "1.7 mm"
$ws.Range("J30").Value = "995.8 hPa"
$ws.Range("E31").Value = "2026-02-13 17:19:36"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "70%"
$ws.Range("I31").Value = "0.8 mm"
$ws.Range("J31").Value = "994.8 hPa"
$ws.Range("E32").Value = "2026-02-13 17:19:38"
$ws.Range("I32").Value = "22.9 mm"
$ws.Range("L32").Value = "38.9 km/h - 283º 16:57 TU"
$ws.Range("E33").Value = "2026-02-13 17:19:41"
$ws.Range("J33").Value = "997.8 hPa"
$ws.Range("K33").Value = "2.7 MJ/m2"
$ws.Range("E34").Value = "2026-02-13 17:19:43"
$ws.Range("G34").Value = "103 cm"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "76%"
$ws.Range("I34").Value = "9.1 mm"
$ws.Range("E35").Value = "2026-02-13 17:19:46"
$ws.Range("I35").Value = "4.7 mm"
$ws.Range("J35").Value = "996.4 hPa"
$ws.Range("E36").Value = "2026-02-13 17:19:48"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "75%"
$ws.Range("I36").Value = "7.3 mm"
$ws.Range("J36").Value = "996.0 hPa"
$ws.Range("E37").Value = "2026-02-13 17:19:51"
$ws.Range("I37").Value = "11.7 mm"
$ws.Range("J37").Value = "997.9 hPa"
$ws.Range("E38").Value = "2026-02-13 17:19:53"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "75%"
$ws.Range("I38").Value = "11.1 mm"
$ws.Range("O38").Value = "9.8 °C"
$ws.Range("E39").Value = "2026-02-13 17:19:56"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "77%"
$ws.Range("I39").Value = "15.7 mm"
$ws.Range("K39").Value = "7.0 MJ/m2"
$ws.Range("E40").Value = "2026-02-13 17:19:58"
$ws.Range("J40").Value = "999.3 hPa"
$ws.Range("O40").Value = "1.4 °C"
$ws.Range("E41").Value = "2026-02-13 17:20:01"
$ws.Range("J41").Value = "996.0 hPa"
$ws.Range("K41").Value = "1.6 MJ/m2"
$ws.Range("L41").Value = "25.9 km/h - 272º 16:42 TU"
$ws.Range("E42").Value = "2026-02-13 17:20:03"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "83%"
$ws.Range("I42").Value = "7.6 mm"
$ws.Range("E43").Value = "2026-02-13 17:20:06"
$ws.Range("I43").Value = "11.4 mm"
$ws.Range("E44").Value = "2026-02-13 17:20:08"
$ws.Range("I44").Value = "3.4 mm"
$ws.Range("K44").Value = "5.6 MJ/m2"
$ws.Range("O44").Value = "-4.1 °C"
$ws.Range("E45").Value = "2026-02-13 17:20:10"
$ws.Range("I45").Value = "0.1 mm"
$ws.Range("J45").Value = "994.1 hPa"
$ws.Range("E46").Value = "2026-02-13 17:20:13"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "91%"
$ws.Range("I46").Value = "9.1 mm"
$ws.Range("J46").Value = "996.5 hPa"
